$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D12")
$sortRange = $ws.Range("A2")

$rng.Sort($sortRange, 1)
